$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "26.582.85" that Excel would
# otherwise auto-parse into numbers on assignment. Force the range to
# Text format while writing, then restore the original (default) style
# so no lasting formatting change is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.592.14'
$ws.Range("D3").Value = '1.818.43'
$ws.Range("D4").Value = '1.008'
$ws.Range("D5").Value = '308.94'
$ws.Range("D7").Value = '0.4568'
$ws.Range("D8").Value = '0.3670'
$ws.Range("D9").Value = '0.07158'
$ws.Range("D10").Value = '0.8795'
$ws.Range("D11").Value = '0.07785'
$ws.Range("D12").Value = '19.41'
$ws.Range("D13").Value = '1.804.73'
$ws.Range("D15").Value = '6.379'
$ws.Range("D16").Value = '86.42'
$ws.Range("D17").Value = '1.009'
$ws.Range("D18").Value = '0.000008626'
$ws.Range("D20").Value = '26.660.78'
$ws.Range("D21").Value = '14.28'
$ws.Range("D22").Value = '5.006'
$ws.Range("D25").Value = '151.62'
$ws.Range("D26").Value = '17.98'
$ws.Range("D28").Value = '113.15'
$ws.Range("D29").Value = '4.870'
$ws.Range("D30").Value = '0.08700'
$ws.Range("D31").Value = '3.069'
$ws.Range("D32").Value = '4.535'
$ws.Range("D33").Value = '0.7377'
$ws.Range("D34").Value = '2.720'
$ws.Range("D38").Value = '0.01947'
$ws.Range("D39").Value = '0.05127'
$ws.Range("D40").Value = '2.905'
$ws.Range("D41").Value = '7.025'
$ws.Range("D42").Value = '0.5035'
$ws.Range("D43").Value = '0.1561'
$ws.Range("D44").Value = '8.209'
$ws.Range("D45").Value = '1.008'
$ws.Range("D46").Value = '0.4635'
$ws.Range("D47").Value = '10.04'
$ws.Range("D48").Value = '101.19'
$ws.Range("D49").Value = '1.597'
$ws.Range("D50").Value = '0.06010'
$ws.Range("D51").Value = '64.68'

# Restore the default (no explicit number format) style on column D.
$ws.Range("D2:D51").Style = "Normal"

# Column E holds percentage-change text (with surrounding spaces), which
# Excel keeps as plain text already, so no special handling is required.
$ws.Range("E2").Value = '  -2.36%  '
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("E16").Value = '  -5.21%  '
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("E18").Value = '  -3.40%  '
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("E27").Value = '  +1.58%  '
$ws.Range("E28").Value = '  -2.42%  '
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("E31").Value = '  -2.21%  '
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("E33").Value = '  -3.93%  '
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("E43").Value = '  -4.18%  '
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("E48").Value = '  -1.46%  '
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("E51").Value = '  -1.02%  '
